$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 7, shifting existing rows 7+ (incl. the
# 1990-2019 FxE year rows) down by one.
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with the new "output" / "configuration_fxe"
# FxE matrix entry (nuclear, enabled).
$ws.Cells.Item(7, 1).Value = "CHE"
$ws.Cells.Item(7, 2).Value = "ext_nuclear"
$ws.Cells.Item(7, 3).Value = "output"
$ws.Cells.Item(7, 4).Value = "configuration_fxe"
$ws.Cells.Item(7, 6).Value = "nuclear"
$ws.Cells.Item(7, 7).Value = 1

# The row insert re-serializes the shifted-down cells; re-assert G16 (the old
# G15 "2.1" cost value) so it doesn't pick up floating-point noise.
$ws.Cells.Item(16, 7).Value = 2.1

# The autofilter range grew by one row (table now spans through row 573);
# refresh it explicitly and sync the hidden _FilterDatabase name to match.
$ws.AutoFilterMode = $false
$ws.Range("A5:L573").AutoFilter()

foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$5:`$L`$573"
    }
}

# Update selection to match the post-edit state.
$ws.Range("E7").Select()
